# This script inserts a new data row right before the current row 20
# (pushing the existing row 20..96 down to 21..97, preserving their
# values) and populates the newly created row 20 with a new record for
# "Puerro" at "Vega Central Mapocho de Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20; Excel shifts rows 20:96 down to 21:97.
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new record.
$ws.Cells.Item(20, 1).Value  = 9
$ws.Cells.Item(20, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value  = "Metropolitana"
$ws.Cells.Item(20, 4).Value  = 44707
$ws.Cells.Item(20, 5).Value  = 13
$ws.Cells.Item(20, 6).Value  = 100112005
$ws.Cells.Item(20, 7).Value  = "Puerro"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 70
$ws.Cells.Item(20, 11).Value = 6000
$ws.Cells.Item(20, 12).Value = 7000
$ws.Cells.Item(20, 13).Value = 6571
$ws.Cells.Item(20, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(20, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(20, 16).Value = 329
$ws.Cells.Item(20, 17).Value = 20
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# Keep the same date number format used by the other rows in column D.
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
